$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1461.75
$ws.Range("I12").Value = 282.66666
$ws.Range("K12").Value = 282.66666
$ws.Range("M12").Value = -112.66666
$ws.Range("H86").Value = 2400
$ws.Range("J86").Value = 2666.6667
$ws.Range("L86").Value = 2666.6667
$ws.Range("N86").Value = -4912.6667
$ws.Range("H89").Value = 2400
$ws.Range("J89").Value = 2666.6667
$ws.Range("L89").Value = 13333.3335
$ws.Range("N89").Value = -24565.3335
$ws.Range("H100").Value = 1862.7778
$ws.Range("I100").Value = 1551.8
$ws.Range("J100").Value = 2251.5
$ws.Range("K100").Value = 1551.8
$ws.Range("L100").Value = 2251.5
$ws.Range("M100").Value = -1010.8
$ws.Range("N100").Value = -3333.5
$ws.Range("H101").Value = 3005.111
$ws.Range("I101").Value = 517.8182
$ws.Range("K101").Value = 1553.4546
$ws.Range("M101").Value = 68.54539999999997
$ws.Range("H112").Value = 2148.1667
$ws.Range("J112").Value = 2168.647
$ws.Range("L112").Value = 6505.941
$ws.Range("N112").Value = -8721.940999999999
$ws.Range("H132").Value = 1298.9
$ws.Range("I132").Value = 1311
$ws.Range("K132").Value = 3933
$ws.Range("M132").Value = -1403

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 1500
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = -1327
$ws.Range("H19").Value = 689
$ws.Range("I19").Value = 649
$ws.Range("J19").Value = 769
$ws.Range("K19").Value = 649
$ws.Range("L19").Value = 769
$ws.Range("M19").Value = -420
$ws.Range("N19").Value = -1227
$ws.Range("H22").Value = 7305.3335
$ws.Range("I22").Value = 4016
$ws.Range("J22").Value = 8950
$ws.Range("K22").Value = 4016
$ws.Range("L22").Value = 8950
$ws.Range("M22").Value = -3717
$ws.Range("N22").Value = -9548
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1598
$ws.Range("N25").ClearContents()  # was -10804
$ws.Range("H45").Value = 5535.143
$ws.Range("I45").Value = 2285.3
$ws.Range("J45").Value = 8489.546
$ws.Range("K45").Value = 2285.3
$ws.Range("L45").Value = 8489.546
$ws.Range("M45").Value = -1908.3
$ws.Range("N45").Value = -9243.546
$ws.Range("H74").Value = 2225.5957
$ws.Range("J74").Value = 2995.6667
$ws.Range("L74").Value = 2995.6667
$ws.Range("N74").Value = -4743.6667
$ws.Range("H77").Value = 2225.5957
$ws.Range("J77").Value = 2995.6667
$ws.Range("L77").Value = 14978.3335
$ws.Range("N77").Value = -23714.3335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2866.3333
$ws.Range("I22").Value = 2866.3333
$ws.Range("K22").Value = 2866.3333
$ws.Range("M22").Value = -2693.3333
$ws.Range("H134").Value = 2046.8302
$ws.Range("I134").Value = 2089.4285
$ws.Range("J134").Value = 1525
$ws.Range("K134").Value = 6268.2855
$ws.Range("L134").Value = 4575
$ws.Range("M134").Value = -3733.2855
$ws.Range("N134").Value = -9645

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10165.866
$ws.Range("I58").Value = 5071.2856
$ws.Range("K58").Value = 5071.2856
$ws.Range("M58").Value = -4868.2856
$ws.Range("H134").Value = 8572
$ws.Range("I134").Value = 6408.2
$ws.Range("J134").Value = 12899.6
$ws.Range("K134").Value = 19224.6
$ws.Range("L134").Value = 38698.8
$ws.Range("M134").Value = -16689.6
$ws.Range("N134").Value = -43768.8
$ws.Range("H136").Value = 10165.866
$ws.Range("I136").Value = 5071.2856
$ws.Range("K136").Value = 15213.8568
$ws.Range("M136").Value = -12663.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47734196
$ws.Range("I4").Value = 54293044
$ws.Range("K4").Value = 162879132
$ws.Range("M4").Value = -162879020
$ws.Range("H23").Value = 9243.546
$ws.Range("I23").Value = 51
$ws.Range("K23").Value = 153
$ws.Range("M23").Value = 82
$ws.Range("H131").Value = 14707304
$ws.Range("I131").Value = 125000980
$ws.Range("K131").Value = 375002940
$ws.Range("M131").Value = -374997900
$ws.Range("H137").Value = 11522.5
$ws.Range("I137").Value = 496.25
$ws.Range("K137").Value = 1488.75
$ws.Range("M137").Value = 3611.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()  # was -4288
$ws.Range("H28").Value = 6995
$ws.Range("J28").Value = 6995
$ws.Range("L28").Value = 6995
$ws.Range("N28").Value = -7379
$ws.Range("H107").Value = 633.6667
$ws.Range("I107").Value = 768.5
$ws.Range("K107").Value = 768.5
$ws.Range("M107").Value = 1151.5
$ws.Range("H123").Value = 39153.77
$ws.Range("J123").Value = 39153.77
$ws.Range("L123").Value = 39153.77
$ws.Range("N123").Value = -44053.77
$ws.Range("H126").Value = 2560.6365
$ws.Range("I126").Value = 1695.1538
$ws.Range("J126").Value = 3810.7778
$ws.Range("K126").Value = 5085.4614
$ws.Range("L126").Value = 11432.3334
$ws.Range("M126").Value = -2615.4614
$ws.Range("N126").Value = -16372.3334
$ws.Range("H132").Value = 2739.5386
$ws.Range("I132").Value = 2395.2334
$ws.Range("K132").Value = 7185.7002
$ws.Range("M132").Value = -4655.7002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1699.5
$ws.Range("I7").Value = 1824.5
$ws.Range("J7").Value = 1512
$ws.Range("K7").Value = 1824.5
$ws.Range("L7").Value = 1512
$ws.Range("M7").Value = -1712.5
$ws.Range("N7").Value = -1736
$ws.Range("H10").Value = 4000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()  # was -2780
$ws.Range("H22").Value = 1525
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 1525
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -2714
$ws.Range("H31").Value = 5266.3335
$ws.Range("J31").Value = 5266.3335
$ws.Range("L31").Value = 5266.3335
$ws.Range("N31").Value = -5762.3335
$ws.Range("H61").Value = 53536.55
$ws.Range("I61").Value = 65890.06
$ws.Range("K61").Value = 65890.06
$ws.Range("M61").Value = -65688.06
$ws.Range("H113").Value = 53536.55
$ws.Range("I113").Value = 65890.06
$ws.Range("K113").Value = 65890.06
$ws.Range("M113").Value = -63720.06
$ws.Range("H126").Value = 1699.5
$ws.Range("I126").Value = 1824.5
$ws.Range("J126").Value = 1512
$ws.Range("K126").Value = 5473.5
$ws.Range("L126").Value = 4536
$ws.Range("M126").Value = -3003.5
